# Regenerate s_vals data to filter save games.
# For each data row (rows 2-16), recompute TB/d2S/K/IP (columns B:E) using the
# updated ("filtered") aggregation, and recompute the row's "sum" (column G)
# as the total of those four components. Dates (A) and Win (F) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    ,@(2,  3.272327238179451,  1.626987699542094,  0.1496068669990043,  0.5333859586016987)
    ,@(3,  3.272327238179451,  1.626987699542094,  0.1496068669990043,  0.5333859586016987)
    ,@(4,  1.445647641019636,  1.626987699542094,  0.1496068669990043,  0.5333859586016987)
    ,@(5,  3.272327238179451,  1.626987699542094,  3.223369029078222,   0.5333859586016987)
    ,@(6,  0.6545652718822623, 1.626987699542094,  0.1496068669990043,  0.5333859586016987)
    ,@(7,  1.445647641019636,  1.626987699542094,  0.7210945179870265,  13.86384647080068)
    ,@(8,  3.272327238179451,  1.626987699542094,  3.223369029078222,   0.5333859586016987)
    ,@(9,  3.272327238179451,  1.626987699542094,  3.223369029078222,   0.5333859586016987)
    ,@(10, 3.272327238179451,  1.626987699542094,  0.7210945179870265,  0.5333859586016987)
    ,@(11, 3.272327238179451,  1.626987699542094,  3.223369029078222,   0.5333859586016987)
    ,@(12, 3.272327238179451,  1.626987699542094,  3.223369029078222,   0.5333859586016987)
    ,@(13, 3.272327238179451,  1.626987699542094,  0.7210945179870265,  0.5333859586016987)
    ,@(14, 1.445647641019636,  0.3048912486333797, 0.1496068669990043,  0.5333859586016987)
    ,@(15, 1.445647641019636,  1.626987699542094,  3.223369029078222,   13.86384647080068)
    ,@(16, 0.1169995834814548, 0.04103571897497393,18.71679738969934,   0.5333859586016987)
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $rowData = $newData[$i]
    $r   = $rowData[0]
    $tb  = $rowData[1]
    $d2s = $rowData[2]
    $k   = $rowData[3]
    $ip  = $rowData[4]

    $ws.Cells.Item($r, 2).Value = $tb
    $ws.Cells.Item($r, 3).Value = $d2s
    $ws.Cells.Item($r, 4).Value = $k
    $ws.Cells.Item($r, 5).Value = $ip
    $ws.Cells.Item($r, 7).Value = ($tb + $d2s + $k + $ip)
}
